$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 193
$ws.Range("F5").Value = 174
$ws.Range("F6").Value = 824
$ws.Range("F7").Value = 4228
$ws.Range("F8").Value = 4228
$ws.Range("F12").Value = 6145
$ws.Range("F15").Value = 2352
$ws.Range("F18").Value = 484
$ws.Range("F19").Value = 9258
$ws.Range("F21").Value = 2494
$ws.Range("F22").Value = 195
$ws.Range("F23").Value = 2324
$ws.Range("F24").Value = 2468
$ws.Range("F25").Value = 1400
$ws.Range("F27").Value = 1979
$ws.Range("F30").Value = 332
$ws.Range("F35").Value = 72
$ws.Range("F38").Value = 1221
$ws.Range("F41").Value = 244
$ws.Range("F42").Value = 1556
$ws.Range("F43").Value = 2557
$ws.Range("F45").Value = 930
$ws.Range("F46").Value = 308
$ws.Range("F48").Value = 27

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 1
$ws.Range("F9").Value = 11
$ws.Range("F22").Value = 83

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 905

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 193
$ws.Range("F4").Value = 905
$ws.Range("F11").Value = 174
$ws.Range("F12").Value = 824
$ws.Range("F13").Value = 4228
$ws.Range("F17").Value = 6145
$ws.Range("F21").Value = 484
$ws.Range("F22").Value = 9258
$ws.Range("F24").Value = 2494
$ws.Range("F25").Value = 195
$ws.Range("F26").Value = 2324
$ws.Range("F27").Value = 2468
$ws.Range("F28").Value = 1400
$ws.Range("F30").Value = 1979
$ws.Range("F33").Value = 332
$ws.Range("F36").Value = 72
$ws.Range("F38").Value = 1221
$ws.Range("F41").Value = 244
$ws.Range("F42").Value = 1556
$ws.Range("F43").Value = 2557
$ws.Range("F44").Value = 930
$ws.Range("F45").Value = 308
$ws.Range("F50").Value = 83
$ws.Range("F51").Value = 83

